$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.333379983901978
$ws.Range("B1").Value = 2.455162048339844
$ws.Range("C1").Value = 4.783269882202148
$ws.Range("D1").Value = 2.498782634735107
$ws.Range("E1").Value = 0.9357627630233765
